# Applies the "Changed the proposed change for 3rd paragraph" edit to
# "Reccomendations for the Game Engine.docx".
#
# The 4th paragraph of the document (the one beginning "We propose to add
# a accessor method to the ActorLocations class...") is rewritten to
# propose a dedicated Locations-tracking class instead of a bare accessor
# method, and the trailing "_GoBack" bookmark is relocated into the
# middle of the rewritten text (where the author's cursor ended up after
# their last edit).
#
# Each Find/Replace below is deliberately scoped to just the surrounding
# "filler" text so that the untouched ActorLocations / GameMap runs (and
# their spell-check <w:proofErr> wrapping) are left exactly as they were.

$d = $word.ActiveDocument

# 1) "We propose to add a accessor method to the " -> new opening sentence.
$rng = $d.Content
$ok1 = $rng.Find.Execute(
    "We propose to add a accessor method to the ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "We propose to add a class solely to govern the Locations of all Actors. From this, we can inherit this class to be used in the game package freely as ",
    2)

# (unchanged run in between: "ActorLocations")

# 2) " class. It would make it easier to find Actors and their Locations
#    outside of the inherited class. Furthermore, ... Location in the "
#    -> new middle sentences, still ending on the same lead-in to GameMap.
$rng = $d.Content
$ok2 = $rng.Find.Execute(
    " class. It would make it easier to find Actors and their Locations outside of the inherited class. Furthermore, it would reduce duplicated code, as we would not need to repeat the process of looping through every Location in the ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " is a protected attribute. It would make it easier to access the ActorLocations without causing any privacy leaks. Furthermore, it would reduce duplicated code, as we would not need to repeat the process of looping through every Location in the ",
    2)

# (unchanged run in between: "GameMap")

# 3) " to find the Actors. However, it would increase privacy leaks as the "
#    -> keep the first clause, change the consequence being described.
$rng = $d.Content
$ok3 = $rng.Find.Execute(
    " to find the Actors. However, it would increase privacy leaks as the ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " to find the Actors. However, it would increase dependencies as it would mean that this new class depends on the GameMap to update its ",
    2)

# (unchanged run in between: "ActorLocations")

# 4) " method can be accessed anywhere instead." -> closing sentence about
#    the encapsulation-boundary design principle.
$rng = $d.Content
$ok4 = $rng.Find.Execute(
    " method can be accessed anywhere instead.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ". Besides that, it would go against the design principle of grouping elements that must depend on each other together inside the encapsulation boundary of a class. ",
    2)

# 5) Relocate the "_GoBack" bookmark (Word's "last edit position" marker)
#    to sit right before "the GameMap to update its ..." -- i.e. where the
#    author's edit in this paragraph ended. Bookmarks.Add with the same
#    name moves an existing bookmark rather than duplicating it, so this
#    also removes it from its old spot at the end of the document.
$bmRng = $d.Content
$okBm = $bmRng.Find.Execute("the GameMap to update its")
$bmPoint = $d.Range($bmRng.Start, $bmRng.Start)
$d.Bookmarks.Add("_GoBack", $bmPoint)

Write-Output "Find results: $ok1 $ok2 $ok3 $ok4 $okBm"
